{"js": "const body = context.document.body;\n\n// Snapshot the paragraph texts once up front so we know exactly which\n// paragraphs are \"Course Description:\" labels and how many there are.\n// (We can't rely on an index computed lazily later because merging\n// paragraphs below shifts every subsequent paragraph's index down by one.)\nconst scanParas = body.paragraphs;\nscanParas.load(\"items/text\");\nawait context.sync();\n\nconst targetOriginalIndices = [];\nfor (let i = 0; i < scanParas.items.length; i++) {\n  if (scanParas.items[i].text === \"Course Description:\") {\n    targetOriginalIndices.push(i);\n  }\n}\n\nlet offset = 0; // how many paragraphs have been removed (merged away) so far\n\nfor (const originalIndex of targetOriginalIndices) {\n  const idx = originalIndex - offset;\n\n  // Re-fetch paragraphs fresh (indices shift as we merge paragraphs away).\n  const paras = body.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n\n  const labelPara = paras.items[idx];       // \"Course Description:\" paragraph\n  const textPara = paras.items[idx + 1];    // the description-text paragraph right after it\n\n  // Step 1: insert a manual line-break character at the very start of the\n  // description-text paragraph, while the two paragraphs are still\n  // separate. Doing it now (rather than after merging) makes the break\n  // land as the leading content of its own run, matching\n  // `<w:r><w:br/><w:t>...</w:t></w:r>`.\n  const startOfTextPara = textPara.getRange(\"Start\");\n  startOfTextPara.insertText(\"\\u000b\", \"Before\");\n  await context.sync();\n\n  // Step 2: bold just the \"Course Description:\" label text (search gives a\n  // tight range that excludes the paragraph mark, so the paragraph-mark\n  // run properties are left untouched).\n  const boldResults = labelPara.search(\"Course Description:\");\n  boldResults.load(\"items\");\n  await context.sync();\n  boldResults.items[0].font.bold = true;\n  await context.sync();\n\n  // Step 3: merge the label paragraph with the (now break-prefixed)\n  // description-text paragraph by deleting the paragraph mark between\n  // them.\n  const mergeParas = body.paragraphs;\n  mergeParas.load(\"items\");\n  await context.sync();\n  const labelParaAgain = mergeParas.items[idx];\n  const textParaAgain = mergeParas.items[idx + 1];\n  const endOfLabel = labelParaAgain.getRange(\"End\");\n  const startOfTextAgain = textParaAgain.getRange(\"Start\");\n  const paraMarkRange = endOfLabel.expandTo(startOfTextAgain);\n  paraMarkRange.delete();\n  await context.sync();\n\n  // Step 4: apply the \"List Bullet\" paragraph style to the merged\n  // paragraph. This must happen after the merge, because merging adopts\n  // the trailing paragraph's (style-less) paragraph properties.\n  const styleParas = body.paragraphs;\n  styleParas.load(\"items\");\n  await context.sync();\n  const mergedPara = styleParas.items[idx];\n  mergedPara.style = \"List Bullet\";\n  await context.sync();\n\n  offset += 1;\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Walk every paragraph; whenever we find one whose entire text is exactly\n# \"Course Description:\" (followed by the paragraph mark), merge it with the\n# paragraph that immediately follows it:\n#   - turn the paragraph-mark that separates the two paragraphs into a\n#     manual line break (w:br) so the description text becomes part of the\n#     same paragraph,\n#   - bold just the \"Course Description:\" label,\n#   - apply the built-in \"List Bullet\" (ListBullet) paragraph style.\n$i = 1\nwhile ($i -le $d.Paragraphs.Count) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text\n  if ($t -eq \"Course Description:`r\") {\n    $descStart = $p.Range.Start\n    $descEnd = $p.Range.End - 1   # exclude the paragraph mark itself\n\n    # Insert a line-break character at the very start of the following\n    # paragraph (while the two paragraphs are still separate) so it lands\n    # in its own run, ahead of the description text.\n    $pNext = $d.Paragraphs.Item($i + 1)\n    $insRange = $d.Range($pNext.Range.Start, $pNext.Range.Start)\n    $insRange.InsertBefore([string][char]11)\n\n    # Delete the paragraph mark that used to separate the two paragraphs,\n    # merging them into one.\n    $pAgain = $d.Paragraphs.Item($i)\n    $markRange = $d.Range($pAgain.Range.End - 1, $pAgain.Range.End)\n    $markRange.Delete()\n\n    # Bold only the \"Course Description:\" label text.\n    $boldRange = $d.Range($descStart, $descEnd)\n    $boldRange.Font.Bold = $true\n\n    # Style the (now merged) paragraph as a bulleted list item.\n    $pStyled = $d.Paragraphs.Item($i)\n    $pStyled.Style = \"ListBullet\"\n  }\n  $i = $i + 1\n}\n"}
